$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = -0.1349881739519865;  C = 1.932678250060537;  D = 16.79600690525843;  E = 4.09829316975475;   F = 4.19246076256923;   G = 22 }
    3  = @{ B = 0.1207090144436769;   C = 1.786104848928261;  D = 11.46089950717241;  E = 3.385395029708116;  F = 3.466791785362177;  G = 21 }
    4  = @{ B = -0.5555707582142663;  C = 1.050109602367337;  D = 4.446451018744308;  E = 2.108660953957347;  F = 2.087000403843689;  G = 20 }
    5  = @{ B = 0.07663430359571852; C = 0.7265601984156437; D = 1.828123414480074;  E = 1.352081141973393;  F = 1.386898245123047;  G = 19 }
    6  = @{ B = 0.04069694792458084; C = 0.747007170374518;  D = 1.696764412358688;  E = 1.302599098862996;  F = 1.339709078915569;  G = 18 }
    7  = @{ B = -0.003405844617402055; C = 0.6048043842844364; D = 0.6518033129846593; E = 0.8073433674618621; F = 0.8321830900003587; G = 17 }
    8  = @{ B = 0.08389527245345252; C = 0.588137197430979;  D = 0.6090342934726032; E = 0.7804064924592844; F = 0.8013294797905551; G = 16 }
    9  = @{ B = 0.225243323359858;   C = 0.4988262046674656; D = 0.4060786597870128; E = 0.6372430147024075; F = 0.6170298427624463; G = 15 }
    10 = @{ B = 0.2001121816970861;  C = 0.4652087079175589; D = 0.3905605582328028; E = 0.6249484444598633; F = 0.614392722160502;  G = 14 }
    11 = @{ B = 0.2336685148374089;  C = 0.3914692457684976; D = 0.1986247317491;    E = 0.4456733464647622; F = 0.395000932487943;  G = 13 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
